$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (week number + date range) ---
$ws.Range("A8").Value = "Volume 31   Number  6"
$ws.Range("C9").Value = "Report Covering the Week  2/5/2024  Through  2/11/2024"

# --- Cell value / type updates ---
$ws.Range("C14").Copy($ws.Range("C15"))
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 4
$ws.Range("H16").Value = -63.636363636363
$ws.Range("I16").Value = 12
$ws.Range("J16").Value = 14
$ws.Range("K16").Value = -14.285714285714
$ws.Range("L16").Value = -33.333333333333
$ws.Range("M16").Value = 9.090909090909
$ws.Range("N16").Value = -83.333333333333
$ws.Range("C17").Value = 5
$ws.Range("E17").Value = 400
$ws.Range("F17").Value = 8
$ws.Range("G17").Value = 8
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 9
$ws.Range("J17").Value = 12
$ws.Range("K17").Value = -25
$ws.Range("L17").Value = -30.76923076923
$ws.Range("M17").Value = -47.058823529411
$ws.Range("N17").Value = -65.384615384615
$ws.Range("C14").Copy($ws.Range("C18"))
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 4
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = -63.636363636363
$ws.Range("J18").Value = 15
$ws.Range("K18").Value = -26.666666666666
$ws.Range("L18").Value = 10
$ws.Range("N18").Value = -73.809523809523
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -7.692307692307
$ws.Range("G19").Value = 52
$ws.Range("H19").Value = -17.307692307692
$ws.Range("I19").Value = 65
$ws.Range("J19").Value = 73
$ws.Range("K19").Value = -10.958904109589
$ws.Range("L19").Value = -16.666666666666
$ws.Range("M19").Value = -18.75
$ws.Range("N19").Value = -29.347826086956
$ws.Range("C20").Value = 1
$ws.Range("D14").Copy($ws.Range("D20"))
$ws.Range("E14").Copy($ws.Range("E20"))
$ws.Range("F20").Value = 6
$ws.Range("H20").Value = 200
$ws.Range("I20").Value = 9
$ws.Range("K20").Value = 200
$ws.Range("L20").Value = -10
$ws.Range("M20").Value = 350
$ws.Range("N20").Value = -84.745762711864
$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 18
$ws.Range("E21").Value = 5.555555555555
$ws.Range("F21").Value = 67
$ws.Range("G21").Value = 84
$ws.Range("H21").Value = -20.238095238095
$ws.Range("I21").Value = 109
$ws.Range("J21").Value = 117
$ws.Range("K21").Value = -6.837606837606
$ws.Range("L21").Value = -15.503875968992
$ws.Range("M21").Value = -12.096774193548
$ws.Range("N21").Value = -62.925170068027
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = -25
$ws.Range("I22").Value = 3
$ws.Range("J22").Value = 5
$ws.Range("K22").Value = -40
$ws.Range("L22").Value = -40
$ws.Range("M22").Value = 200
$ws.Range("C23").Value = 2
$ws.Range("F23").Value = 4
$ws.Range("H23").Value = 300
$ws.Range("I23").Value = 5
$ws.Range("K23").Value = 66.666666666666
$ws.Range("L23").Value = -44.444444444444
$ws.Range("M23").Value = 66.666666666666
$ws.Range("C24").Value = 13
$ws.Range("D24").Value = 11
$ws.Range("E24").Value = 18.181818181818
$ws.Range("F24").Value = 44
$ws.Range("G24").Value = 44
$ws.Range("I24").Value = 71
$ws.Range("J24").Value = 66
$ws.Range("K24").Value = 7.575757575757
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -20.224719101123
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = -25
$ws.Range("F25").Value = 26
$ws.Range("G25").Value = 29
$ws.Range("H25").Value = -10.344827586206
$ws.Range("I25").Value = 42
$ws.Range("J25").Value = 43
$ws.Range("K25").Value = -2.325581395348
$ws.Range("L25").Value = 110
$ws.Range("M25").Value = 27.272727272727
$ws.Range("C14").Copy($ws.Range("C26"))
$ws.Range("C14").Copy($ws.Range("C27"))
$ws.Range("D14").Copy($ws.Range("D27"))
$ws.Range("E14").Copy($ws.Range("E27"))
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 66.666666666666
$ws.Range("L27").Value = 20
$ws.Range("D16").Copy($ws.Range("D30"))
$ws.Range("D30").Value = 3
$ws.Range("E16").Copy($ws.Range("E30"))
$ws.Range("E30").Value = -100
$ws.Range("G16").Copy($ws.Range("G30"))
$ws.Range("G30").Value = 3
$ws.Range("H16").Copy($ws.Range("H30"))
$ws.Range("H30").Value = -100
$ws.Range("J16").Copy($ws.Range("J30"))
$ws.Range("J30").Value = 3
$ws.Range("K16").Copy($ws.Range("K30"))
$ws.Range("K30").Value = -100
